{"js": "// Replace the date and each two-digit multiplication answer in the document.\n// The mapping below is derived from the diff: every <w:t> run that changes\n// has a unique \"before\" string in this document, so an exact-text,\n// match-case search/replace is safe and unambiguous.\nconst replacements = [\n  [\"2025-10-19 Sunday\", \"2025-10-20 Monday\"],\n  [\"63\u00d784=5292\", \"24\u00d724=576\"],\n  [\"92\u00d791=8372\", \"44\u00d787=3828\"],\n  [\"19\u00d794=1786\", \"85\u00d776=6460\"],\n  [\"21\u00d764=1344\", \"61\u00d745=2745\"],\n  [\"56\u00d757=3192\", \"26\u00d758=1508\"],\n  [\"99\u00d715=1485\", \"26\u00d751=1326\"],\n  [\"63\u00d782=5166\", \"45\u00d731=1395\"],\n  [\"37\u00d791=3367\", \"40\u00d799=3960\"],\n  [\"12\u00d720=240\", \"35\u00d717=595\"],\n  [\"49\u00d758=2842\", \"11\u00d732=352\"],\n  [\"75\u00d743=3225\", \"36\u00d783=2988\"],\n  [\"81\u00d751=4131\", \"24\u00d749=1176\"],\n  [\"77\u00d759=4543\", \"58\u00d792=5336\"],\n  [\"64\u00d714=896\", \"49\u00d760=2940\"],\n  [\"37\u00d794=3478\", \"85\u00d757=4845\"],\n  [\"86\u00d760=5160\", \"15\u00d719=285\"],\n  [\"26\u00d780=2080\", \"35\u00d722=770\"],\n  [\"62\u00d728=1736\", \"99\u00d795=9405\"],\n  [\"91\u00d796=8736\", \"54\u00d797=5238\"],\n  [\"28\u00d764=1792\", \"51\u00d794=4794\"],\n  [\"34\u00d755=1870\", \"41\u00d746=1886\"],\n  [\"61\u00d719=1159\", \"50\u00d792=4600\"],\n  [\"68\u00d713=884\", \"80\u00d774=5920\"],\n  [\"80\u00d785=6800\", \"29\u00d739=1131\"],\n  [\"24\u00d745=1080\", \"91\u00d797=8827\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit multiplication answer.\n# Every \"before\" string below is unique within the document, so a\n# whole-document Find/Replace (ReplaceAll) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-19 Sunday\", \"2025-10-20 Monday\"),\n    @(\"63\u00d784=5292\", \"24\u00d724=576\"),\n    @(\"92\u00d791=8372\", \"44\u00d787=3828\"),\n    @(\"19\u00d794=1786\", \"85\u00d776=6460\"),\n    @(\"21\u00d764=1344\", \"61\u00d745=2745\"),\n    @(\"56\u00d757=3192\", \"26\u00d758=1508\"),\n    @(\"99\u00d715=1485\", \"26\u00d751=1326\"),\n    @(\"63\u00d782=5166\", \"45\u00d731=1395\"),\n    @(\"37\u00d791=3367\", \"40\u00d799=3960\"),\n    @(\"12\u00d720=240\", \"35\u00d717=595\"),\n    @(\"49\u00d758=2842\", \"11\u00d732=352\"),\n    @(\"75\u00d743=3225\", \"36\u00d783=2988\"),\n    @(\"81\u00d751=4131\", \"24\u00d749=1176\"),\n    @(\"77\u00d759=4543\", \"58\u00d792=5336\"),\n    @(\"64\u00d714=896\", \"49\u00d760=2940\"),\n    @(\"37\u00d794=3478\", \"85\u00d757=4845\"),\n    @(\"86\u00d760=5160\", \"15\u00d719=285\"),\n    @(\"26\u00d780=2080\", \"35\u00d722=770\"),\n    @(\"62\u00d728=1736\", \"99\u00d795=9405\"),\n    @(\"91\u00d796=8736\", \"54\u00d797=5238\"),\n    @(\"28\u00d764=1792\", \"51\u00d794=4794\"),\n    @(\"34\u00d755=1870\", \"41\u00d746=1886\"),\n    @(\"61\u00d719=1159\", \"50\u00d792=4600\"),\n    @(\"68\u00d713=884\", \"80\u00d774=5920\"),\n    @(\"80\u00d785=6800\", \"29\u00d739=1131\"),\n    @(\"24\u00d745=1080\", \"91\u00d797=8827\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
